$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")

# Append extra detail text to the two suggestion cells, matching the
# author's edit (extending two existing ideas with additional remarks).
$ws.Range("B8").Value = $ws.Range("B8").Text + ", модульный интерфейс"
$ws.Range("B3").Value = $ws.Range("B3").Text + ", обучение новичков"

# Reflect the scroll/selection state captured in the saved file: the
# window had been scrolled down and the user had clicked near B14.
$ws.Activate()
$ws.Range("A10").Select()
$excel.ActiveWindow.ScrollRow = 10
$ws.Range("B14").Select()
